$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same format (style) used by the existing data rows (column A,
# style index 4 in the original file) to the whole new column U range
# (header row handled separately below), so every row lands on style "4"
# like the target workbook - not the alternating 4/7 banding column A uses.
$ws.Range("A2").Copy()
$ws.Range("U2:U82").PasteSpecial(-4122)

# Header cell U1: copy the format from T1 (previously the last header
# column) so U1 picks up the same header styling.
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$ws.Range("U1").Value = "POSH Win"

# New "POSH Win" flag column: 0 for every non-POSH-core winner row,
# 1 for the rows where POSH-core won (rows 75-82).
for ($row = 2; $row -le 74; $row++) {
    $ws.Range("U$row").Value = 0
}
for ($row = 75; $row -le 82; $row++) {
    $ws.Range("U$row").Value = 1
}

# Restore the clipboard / leave selection where the author left it when
# they saved (scrolled down, landed on Y8).
$ws.Range("Y8").Select()
